# Update TPM-derived statistics in the Efna5-Ephb1 LR-pair sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 0.1728506666666667
$ws.Range("H2").Value = 0.518552
$ws.Range("I2").Value = 0.0840503369699626
$ws.Range("J2").Value = 0.0840503369699626
$ws.Range("M2").Value = 2.718682666666667
$ws.Range("N2").Value = 8.156048
$ws.Range("O2").Value = 0.5434637507613679
$ws.Range("P2").Value = 0.5434637507613679
$ws.Range("Q2").Value = 0.4699261113884445
$ws.Range("R2").Value = 4.229335002496001
$ws.Range("S2").Value = 0.04567831138245274
$ws.Range("T2").Value = 0.04567831138245274

# Row 3 (ECs -> MuSCs)
$ws.Range("G3").Value = 0.1728506666666667
$ws.Range("H3").Value = 0.518552
$ws.Range("I3").Value = 0.0840503369699626
$ws.Range("J3").Value = 0.0840503369699626
$ws.Range("M3").Value = 2.283827
$ws.Range("N3").Value = 6.851481
$ws.Range("O3").Value = 0.4565362492386322
$ws.Range("P3").Value = 0.4565362492386321
$ws.Range("Q3").Value = 0.3947610195013334
$ws.Range("R3").Value = 3.552849175512
$ws.Range("S3").Value = 0.03837202558750986
$ws.Range("T3").Value = 0.03837202558750986

# Row 4 (FAPs -> ECs)
$ws.Range("I4").Value = 0.6650661694281633
$ws.Range("J4").Value = 0.6650661694281633
$ws.Range("M4").Value = 2.718682666666667
$ws.Range("N4").Value = 8.156048
$ws.Range("O4").Value = 0.5434637507613679
$ws.Range("P4").Value = 0.5434637507613679
$ws.Range("Q4").Value = 3.718390313260444
$ws.Range("R4").Value = 33.465512819344
$ws.Range("S4").Value = 0.361439354941925
$ws.Range("T4").Value = 0.361439354941925

# Row 5 (FAPs -> MuSCs)
$ws.Range("I5").Value = 0.6650661694281633
$ws.Range("J5").Value = 0.6650661694281633
$ws.Range("M5").Value = 2.283827
$ws.Range("N5").Value = 6.851481
$ws.Range("O5").Value = 0.4565362492386322
$ws.Range("P5").Value = 0.4565362492386321
$ws.Range("Q5").Value = 3.123630535510333
$ws.Range("R5").Value = 28.112674819593
$ws.Range("S5").Value = 0.3036268144862384
$ws.Range("T5").Value = 0.3036268144862383

# Row 6 (MuSCs -> ECs)
$ws.Range("G6").Value = 0.5159453333333334
$ws.Range("H6").Value = 1.547836
$ws.Range("I6").Value = 0.2508834936018741
$ws.Range("J6").Value = 0.2508834936018741
$ws.Range("M6").Value = 2.718682666666667
$ws.Range("N6").Value = 8.156048
$ws.Range("O6").Value = 0.5434637507613679
$ws.Range("P6").Value = 0.5434637507613679
$ws.Range("Q6").Value = 1.402691634680889
$ws.Range("R6").Value = 12.624224712128
$ws.Range("S6").Value = 0.1363460844369901
$ws.Range("T6").Value = 0.1363460844369901

# Row 7 (MuSCs -> MuSCs)
$ws.Range("G7").Value = 0.5159453333333334
$ws.Range("H7").Value = 1.547836
$ws.Range("I7").Value = 0.2508834936018741
$ws.Range("J7").Value = 0.2508834936018741
$ws.Range("M7").Value = 2.283827
$ws.Range("N7").Value = 6.851481
$ws.Range("O7").Value = 0.4565362492386322
$ws.Range("P7").Value = 0.4565362492386321
$ws.Range("Q7").Value = 1.178329882790667
$ws.Range("R7").Value = 10.604968945116
$ws.Range("S7").Value = 0.114537409164884
$ws.Range("T7").Value = 0.114537409164884
